$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.226.04'
$ws.Range("E2").Value = '  -0.79%  '
$ws.Range("D3").Value = '1.807.98'
$ws.Range("E3").Value = '  +1.16%  '
$ws.Range("D4").Value = "'" + '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = "'" + '223.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = "'" + '33.02'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.32%  '
$ws.Range("D9").Value = "'" + '0.288'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.68%  '
$ws.Range("D10").Value = "'" + '0.0718'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.68%  '
$ws.Range("D11").Value = "'" + '0.0928'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.65%  '
$ws.Range("D12").Value = '2.067.86'
$ws.Range("E12").Value = '  +1.18%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.824.73'
$ws.Range("E13").Value = '  +2.45%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = "'" + '11.07'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.87%  '
$ws.Range("D15").Value = "'" + '0.632'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.26%  '
$ws.Range("D16").Value = '34.264.29'
$ws.Range("E16").Value = '  -0.72%  '
$ws.Range("D17").Value = "'" + '4.24'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.97%  '
$ws.Range("D18").Value = "'" + '68.91'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D19").Value = "'" + '247.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.33%  '
$ws.Range("D20").Value = '0.0₃0788'
$ws.Range("E20").Value = '  +0.84%  '
$ws.Range("E21").Value = '  +5.82%  '
$ws.Range("D22").Value = "'" + '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").Value = "'" + '4.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.12%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").Value = "'" + '159.81'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("D26").Value = "'" + '16.62'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.47%  '
$ws.Range("D27").Value = "'" + '7.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.13%  '
$ws.Range("E28").Value = '  -0.87%  '
$ws.Range("D29").Value = "'" + '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  +2.65%  '
$ws.Range("D31").Value = "'" + '3.74'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("E32").Value = '  +1.83%  '
$ws.Range("D33").Value = "'" + '3.54'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.50%  '
$ws.Range("E34").Value = '  -0.72%  '
$ws.Range("D35").Value = '1.421.04'
$ws.Range("E35").Value = '  -1.06%  '
$ws.Range("E36").Value = '  +2.20%  '
$ws.Range("E37").Value = '  +0.55%  '
$ws.Range("E38").Value = '  -1.04%  '
$ws.Range("D39").Value = "'" + '0.947'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.30%  '
$ws.Range("D40").Value = "'" + '80.93'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.88%  '
$ws.Range("E41").Value = '  -2.37%  '
$ws.Range("E42").Value = '  +0.78%  '
$ws.Range("E43").Value = '  +4.17%  '
$ws.Range("D44").Value = "'" + '5.96'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.73%  '
$ws.Range("D45").Value = "'" + '108.36'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.22%  '
$ws.Range("D46").Value = "'" + '0.0497'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.24%  '
$ws.Range("E47").Value = '  -1.15%  '
$ws.Range("D48").Value = '1.966.30'
$ws.Range("E48").Value = '  +1.03%  '
$ws.Range("D49").Value = "'" + '12.08'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.23%  '
$ws.Range("E51").Value = '  +3.02%  '
